$d = $word.ActiveDocument

# Paragraph 1: "ccbdb" -> "Sujet 1"
$d.Paragraphs.Item(1).Range.Text = "Sujet 1"

# Paragraph 2: "dcdcb" -> "cdbbc adaca daccc acacc dcbdc " (becomes the sole
# replacement for the former paragraphs 2-5)
$d.Paragraphs.Item(2).Range.Text = "cdbbc adaca daccc acacc dcbdc "

# Remove the now-redundant former paragraphs 3, 4 and 5 (deleting a
# paragraph's Range removes its text together with its paragraph mark).
$d.Paragraphs.Item(5).Range.Delete()
$d.Paragraphs.Item(4).Range.Delete()
$d.Paragraphs.Item(3).Range.Delete()
